$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1. "...installed and winscp please follow..." -> "...installed and WinSCP please follow..."
#    (the paragraph is later re-split into several runs in the authored diff,
#    purely for styling/bookmark placement purposes; the visible text only
#    changes the capitalisation of "winscp" -> "WinSCP")
# --------------------------------------------------------------------------
$d.Content.Find.Execute(
    "If you do not have the tomcat component installed and winscp please follow these lines",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "If you do not have the tomcat component installed and WinSCP please follow these lines",
    2) | Out-Null

# --------------------------------------------------------------------------
# 2. Replace the "Copy tomcat8.tar, winscp577setup.exe and putty.exe from usb
#    drive to a new directory on your computer." bullet with the new
#    "Install WinSCP on your local machine: <link>" bullet.
# --------------------------------------------------------------------------
$copyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("Copy tomcat8.tar")) {
        $copyPara = $cand
        break
    }
}

$r = $copyPara.Range
$r.End = $r.End - 1
$r.Text = "Install WinSCP on your local machine: https://winscp.net/eng/download.php "

# Turn the URL text into a real hyperlink.
$urlRange = $r.Duplicate
$urlRange.Find.Execute("https://winscp.net/eng/download.php", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($urlRange, "https://winscp.net/eng/download.php") | Out-Null

# --------------------------------------------------------------------------
# 3. Remove the now-obsolete "Open a command prompt..." / docker load / docker
#    run / docker start / "Install winscp on you local machine." bullets (and
#    their surrounding blank lines) - the whole block between the "Open a
#    command prompt..." bullet and the "Open winscp and connect..." bullet.
# --------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    $t = $cand.Range.Text
    if ($startPara -eq $null -and $t.StartsWith("Open a command prompt")) {
        $startPara = $cand
    }
    if ($t.StartsWith("Open winscp and connect to the tomcat server")) {
        $endPara = $cand
        break
    }
}

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
$delRange.Delete()

# --------------------------------------------------------------------------
# 4. Update the wording of the "Open winscp and connect..." bullet.
# --------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Open winscp and connect to the tomcat server on docker like in the printscreen (the username and password are root and root):",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Open WinSCP and connect to the tomcat server on Docker like in the screenshot (the username and password are root and root):",
    2) | Out-Null

Write-Output "done"
